# Improved documentation of ArrayPanel.
#
# 1. Duplicate the existing "ArrayPanel.Orientation" slide so the deck gets
#    a second slide (the duplicate keeps the original, unedited diagram).
# 2. Clean up / simplify the diagram on the original (now first) slide by
#    removing the padding-measurement callout and its helper shapes, then
#    resize & reposition the remaining shapes so the diagram reflows nicely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Duplicate slide 1 -> becomes new slide 2 (untouched copy) ---------
$null = $s.Duplicate()

# --- 2. Remove the padding-measurement shapes from slide 1 ----------------
# "ArrayPanel.Padding" label
$s.Shapes.Item("TextBox 20").Delete()
# Arrow from the padding label to the inner padding gap
$s.Shapes.Item("Straight Arrow Connector 21").Delete()
# Small padding-gap marker rectangles
$s.Shapes.Item("Rectangle 24").Delete()
$s.Shapes.Item("Rectangle 25").Delete()
$s.Shapes.Item("Rectangle 30").Delete()
$s.Shapes.Item("Rectangle 31").Delete()
# Arrow from the padding label to the other padding-gap marker
$s.Shapes.Item("Straight Arrow Connector 34").Delete()

# --- 3. Resize / reposition the surviving shapes ---------------------------

# Outer panel rectangle
$shp = $s.Shapes.Item("Rectangle 3")
$shp.Left = 201.23669291338584
$shp.Top = 171.07464566929133
$shp.Width = 232.44086614173227
$shp.Height = 73.7007874015748

# Baseline arrow under the panel
$shp = $s.Shapes.Item("Straight Arrow Connector 8")
$shp.Left = 201.23669291338584
$shp.Top = 270.0
$shp.Width = 232.44086614173227
$shp.Height = 0.0

# Right brace spanning the panel width
$shp = $s.Shapes.Item("Right Brace 11")
$shp.Left = 306.5516535433071
$shp.Top = -0.07283464566929133
$shp.Width = 21.810944881889764
$shp.Height = 232.44086614173227

# "ArrayPanel.Spacing" label
$shp = $s.Shapes.Item("TextBox 14")
$shp.Left = 342.65574803149605
$shp.Top = 136.65141732283465

# Arrow from spacing label to 3rd child
$shp = $s.Shapes.Item("Straight Arrow Connector 16")
$shp.Left = 359.9767716535433
$shp.Top = 158.4623622047244
$shp.Width = 36.8503937007874
$shp.Height = 49.46267716535433

# Arrow from spacing label to 2nd child
$shp = $s.Shapes.Item("Straight Arrow Connector 18")
$shp.Left = 280.6067716535433
$shp.Top = 158.4623622047244
$shp.Width = 116.2203937007874
$shp.Height = 49.46267716535433

# "Available width distributed equally among child elements" caption
$shp = $s.Shapes.Item("TextBox 37")
$shp.Left = 214.5867716535433
$shp.Top = 306.31566929133857

# Arrow from caption to 1st child
$shp = $s.Shapes.Item("Straight Arrow Connector 39")
$shp.Width = 57.17488188976378
$shp.Height = 61.540236220472444

# Arrow from caption to 2nd child
$shp = $s.Shapes.Item("Straight Arrow Connector 41")
$shp.Left = 295.2620472440945
$shp.Top = 244.77543307086614
$shp.Width = 22.19511811023622
$shp.Height = 61.540236220472444
$shp.HorizontalFlip = $false

# Arrow from caption to 3rd child
$shp = $s.Shapes.Item("Straight Arrow Connector 43")
$shp.Left = 295.2620472440945
$shp.Top = 244.77543307086614
$shp.Width = 101.56511811023623
$shp.Height = 61.540236220472444

# "ArrayPanel.Orientation: Horizontal" caption
$shp = $s.Shapes.Item("TextBox 9")
$shp.Left = 192.25582677165355
$shp.Top = 268.7357480314961
